$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112181582
$ws.Range("B2").Value = 81385
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 1312
$ws.Range("F2").Value = 'Gammelgransskål'
$ws.Range("G2").Value = 'Pseudographis pinicola'
$ws.Range("H2").Value = '(Nyl.) Rehm'
$ws.Range("Q2").Value = 772409
$ws.Range("R2").Value = 7120320

# Row 3
$ws.Range("A3").Value = 112182926
$ws.Range("B3").Value = 5113
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 100526
$ws.Range("F3").Value = 'Bronshjon'
$ws.Range("G3").Value = 'Callidium coriaceum'
$ws.Range("H3").Value = 'Paykull, 1800'
$ws.Range("Q3").Value = 772357
$ws.Range("R3").Value = 7120234
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = '2023-09-06'
$ws.Range("Y3").ClearFormats()
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = '2023-09-06'
$ws.Range("AA3").ClearFormats()

# Row 4
$ws.Range("A4").Value = 112181514
$ws.Range("B4").Value = 89553
$ws.Range("Q4").Value = 772353
$ws.Range("R4").Value = 7120281
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = '2023-09-06'
$ws.Range("Y4").ClearFormats()
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = '2023-09-06'
$ws.Range("AA4").ClearFormats()

# Row 5
$ws.Range("B5").Value = 89553

# Row 6
$ws.Range("A6").Value = 112181500
$ws.Range("B6").Value = 89517
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 5447
$ws.Range("F6").Value = 'Vedticka'
$ws.Range("G6").Value = 'Fuscoporia viticola'
$ws.Range("H6").Value = '(Schwein.) Murrill'
$ws.Range("Q6").Value = 772346
$ws.Range("R6").Value = 7120286

# Row 7
$ws.Range("B7").Value = 89499

# Row 8
$ws.Range("A8").Value = 112181509
$ws.Range("B8").Value = 89553
$ws.Range("E8").Value = 1202
$ws.Range("F8").Value = 'Ullticka'
$ws.Range("G8").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H8").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q8").Value = 772347
$ws.Range("R8").Value = 7120237
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = '2023-09-07'
$ws.Range("Y8").ClearFormats()
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = '2023-09-07'
$ws.Range("AA8").ClearFormats()

# Row 9
$ws.Range("A9").Value = 112181511
$ws.Range("B9").Value = 89553
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 1202
$ws.Range("F9").Value = 'Ullticka'
$ws.Range("G9").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H9").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q9").Value = 772359
$ws.Range("R9").Value = 7120174
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = '2023-09-07'
$ws.Range("Y9").ClearFormats()
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = '2023-09-07'
$ws.Range("AA9").ClearFormats()
